$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "empty target" rows appended below the existing data (A:B only, C left blank).
$newRows = @(
  @(56, "g2 ", "age"),
  @(57, "g1",  "age"),
  @(58, "g1",  "car "),
  @(59, "g1",  "cst"),
  @(60, "g2 ", "cst"),
  @(61, "g1",  "edu"),
  @(62, "g1",  "eth"),
  @(63, "g1",  "gen"),
  @(64, "g2 ", "gen"),
  @(65, "g1 ", "inc"),
  @(66, "g1",  "inv"),
  @(67, "g1",  "occ"),
  @(68, "g1",  "par"),
  @(69, "g1",  "ris"),
  @(70, "g1",  "rso"),
  @(71, "g1",  "rst"),
  @(72, "g2 ", "rst"),
  @(73, "g1",  "ses"),
  @(74, "g1",  "sex"),
  @(75, "g2 ", "sex")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Re-point the AutoFilter at the grown range (toggle off/on so the stored
# range actually refreshes instead of just flipping AutoFilterMode).
$lastRow = 75
$ws.Range("A1:C$lastRow").AutoFilter() | Out-Null
$ws.Range("A1:C$lastRow").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$$lastRow"
    }
}

# Match the author's final selection/scroll target.
$ws.Range("C55").Select() | Out-Null
